# Update TranslatorTestingModel sheets: replace the trailing biolink_*
# qualifier columns with id/name/description/tags on TestAsset,
# AcceptanceTestAsset, and TestEdgeData, shrinking the used range.

$wb = $excel.ActiveWorkbook

# --- TestAsset: columns S1:AA1 -> S1:V1 (id, name, description, tags) ---
$ws = $wb.Worksheets.Item("TestAsset")
$ws.Range("S1:AA1").ClearContents()
$ws.Range("S1").Value = "id"
$ws.Range("T1").Value = "name"
$ws.Range("U1").Value = "description"
$ws.Range("V1").Value = "tags"

# --- AcceptanceTestAsset: columns AC1:AK1 -> AC1:AF1 (id, name, description, tags) ---
$ws = $wb.Worksheets.Item("AcceptanceTestAsset")
$ws.Range("AC1:AK1").ClearContents()
$ws.Range("AC1").Value = "id"
$ws.Range("AD1").Value = "name"
$ws.Range("AE1").Value = "description"
$ws.Range("AF1").Value = "tags"

# --- TestEdgeData: columns S1:AA1 -> S1:V1 (id, name, description, tags) ---
$ws = $wb.Worksheets.Item("TestEdgeData")
$ws.Range("S1:AA1").ClearContents()
$ws.Range("S1").Value = "id"
$ws.Range("T1").Value = "name"
$ws.Range("U1").Value = "description"
$ws.Range("V1").Value = "tags"
